$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price / volume(1h) figures scraped on Tue Apr 30 23:06:49 UTC 2024.
# Column D (Price) values are quote-prefixed so Excel stores them as text
# (several look numeric, e.g. '1.00', '0.110' -- without the prefix Excel's
# COM layer would silently coerce them to numbers and drop the formatting).

$ws.Range("D2").Value = "'60.568.74"
$ws.Range("E2").Value = '  -5.30%  '
$ws.Range("D3").Value = "'3.006.93"
$ws.Range("E3").Value = '  -6.56%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = "'577.58"
$ws.Range("E5").Value = '  -2.83%  '
$ws.Range("D6").Value = "'126.79"
$ws.Range("E6").Value = '  -8.00%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").Value = "'3.002.92"
$ws.Range("E8").Value = '  -6.70%  '
$ws.Range("D9").Value = "'0.499"
$ws.Range("E9").Value = '  -3.24%  '
$ws.Range("E10").Value = '  -7.54%  '
$ws.Range("E11").Value = '  -2.73%  '
$ws.Range("D12").Value = "'0.440"
$ws.Range("E12").Value = '  -3.60%  '
$ws.Range("E13").Value = '  -7.25%  '
$ws.Range("D14").Value = "'32.67"
$ws.Range("E14").Value = '  -6.79%  '
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").Value = "'3.506.50"
$ws.Range("E16").Value = '  -6.40%  '
$ws.Range("D17").Value = "'3.009.25"
$ws.Range("E17").Value = '  -6.43%  '
$ws.Range("D18").Value = "'60.511.10"
$ws.Range("E18").Value = '  -5.30%  '
$ws.Range("D19").Value = "'6.37"
$ws.Range("E19").Value = '  -3.78%  '
$ws.Range("D20").Value = "'432.77"
$ws.Range("E20").Value = '  -7.54%  '
$ws.Range("D21").Value = "'13.11"
$ws.Range("E21").Value = '  -6.78%  '
$ws.Range("D22").Value = "'0.665"
$ws.Range("E22").Value = '  -6.12%  '
$ws.Range("D23").Value = "'7.02"
$ws.Range("E23").Value = '  -9.11%  '
$ws.Range("D24").Value = "'12.86"
$ws.Range("E24").Value = '  -4.99%  '
$ws.Range("D25").Value = "'79.53"
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("E28").Value = '  -5.20%  '
$ws.Range("E29").Value = '  -6.88%  '
$ws.Range("E30").Value = '  -8.18%  '
$ws.Range("D31").Value = "'6.15"
$ws.Range("E31").Value = '  -10.47%  '
$ws.Range("D32").Value = "'25.33"
$ws.Range("E32").Value = '  -8.30%  '
$ws.Range("D33").Value = "'0.0937"
$ws.Range("E33").Value = '  -9.47%  '
$ws.Range("D34").Value = "'2.16"
$ws.Range("E34").Value = '  -11.07%  '
$ws.Range("D35").Value = "'0.959"
$ws.Range("E35").Value = '  -7.82%  '
$ws.Range("D36").Value = "'5.61"
$ws.Range("E36").Value = '  -5.32%  '
$ws.Range("D37").Value = "'50.05"
$ws.Range("E37").Value = '  -3.24%  '
$ws.Range("D38").Value = "'0.0₃0666"
$ws.Range("E38").Value = '  -9.41%  '
$ws.Range("D39").Value = "'8.33"
$ws.Range("E39").Value = '  +2.13%  '
$ws.Range("D40").Value = "'0.0361"
$ws.Range("E40").Value = '  -8.58%  '
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").Value = "'388.01"
$ws.Range("E41").Value = '  -4.97%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = "'0.110"
$ws.Range("E42").Value = '  -2.96%  '
$ws.Range("D43").Value = "'2.52"
$ws.Range("E43").Value = '  -9.11%  '
$ws.Range("D44").Value = "'2.661.32"
$ws.Range("E44").Value = '  -6.12%  '
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").Value = "'0.236"
$ws.Range("E46").Value = '  -7.98%  '
$ws.Range("D47").Value = "'2.02"
$ws.Range("E47").Value = '  -7.04%  '
$ws.Range("D48").Value = "'118.06"
$ws.Range("E48").Value = '  -7.37%  '
$ws.Range("E49").Value = '  -4.43%  '
$ws.Range("D50").Value = "'23.79"
$ws.Range("E50").Value = '  -7.90%  '
$ws.Range("E51").Value = '  +3.46%  '
